# Update LR-pair data with new TPM-based values (Cxcl13-Cxcr3).
# Target clusters expand from {Resolving-Mac, Inflammatory-Mac} to
# {ECs, Inflammatory-Mac, Neutrophils, Resolving-Mac}, and a new sending
# cluster 'MuSCs' is added alongside the existing 'FAPs', producing 8 data rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (updated)
$ws.Range("D2").Value = "ECs"
$ws.Range("G2").Value = 7.767740666666666
$ws.Range("I2").Value = 0.957755623847744
$ws.Range("J2").Value = 0.9714348434930534
$ws.Range("K2").Value = 1
$ws.Range("L2").Value = 0.5
$ws.Range("M2").Value = 0.358189
$ws.Range("N2").Value = 0.716378
$ws.Range("O2").Value = 0.052978466298774
$ws.Range("P2").Value = 0.03595390511601538
$ws.Range("Q2").Value = 2.782319261652666
$ws.Range("R2").Value = 16.693915569916
$ws.Range("S2").Value = 0.05074042404047897
$ws.Range("T2").Value = 0.03492687618934049

# Row 3 (updated)
$ws.Range("D3").Value = "Inflammatory-Mac"
$ws.Range("G3").Value = 7.767740666666666
$ws.Range("I3").Value = 0.957755623847744
$ws.Range("J3").Value = 0.9714348434930534
$ws.Range("K3").Value = 2
$ws.Range("L3").Value = 0.6666666666666666
$ws.Range("M3").Value = 0.574863
$ws.Range("N3").Value = 1.724589
$ws.Range("O3").Value = 0.08502595018806307
$ws.Range("P3").Value = 0.08655445766079338
$ws.Range("Q3").Value = 4.465386702861999
$ws.Range("R3").Value = 40.18848032575799
$ws.Range("S3").Value = 0.08143408196561555
$ws.Range("T3").Value = 0.08408201603133894

# Row 4 (updated)
$ws.Range("A4").Value = "FAPs"
$ws.Range("D4").Value = "Neutrophils"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 7.767740666666666
$ws.Range("H4").Value = 23.303222
$ws.Range("I4").Value = 0.957755623847744
$ws.Range("J4").Value = 0.9714348434930534
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 5.590969
$ws.Range("N4").Value = 16.772907
$ws.Range("O4").Value = 0.8269404217996371
$ws.Range("P4").Value = 0.8418062905306279
$ws.Range("Q4").Value = 43.42919726737266
$ws.Range("R4").Value = 390.862775406354
$ws.Range("S4").Value = 0.792006839565628
$ws.Range("T4").Value = 0.8177599620930884

# Row 5 (updated)
$ws.Range("A5").Value = "FAPs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 7.767740666666666
$ws.Range("H5").Value = 23.303222
$ws.Range("I5").Value = 0.957755623847744
$ws.Range("J5").Value = 0.9714348434930534
$ws.Range("M5").Value = 0.237009
$ws.Range("N5").Value = 0.7110270000000001
$ws.Range("O5").Value = 0.0350551617135259
$ws.Range("P5").Value = 0.03568534669256324
$ws.Range("Q5").Value = 1.841024447666
$ws.Range("R5").Value = 16.569220028994
$ws.Range("S5").Value = 0.03357427827602155
$ws.Range("T5").Value = 0.03466598917928553

# Row 6 (new)
$ws.Range("A6").Value = "MuSCs"
$ws.Range("B6").Value = "Cxcl13"
$ws.Range("C6").Value = "Cxcr3"
$ws.Range("D6").Value = "ECs"
$ws.Range("E6").Value = 1
$ws.Range("F6").Value = 0.5
$ws.Range("G6").Value = 0.342617
$ws.Range("H6").Value = 0.685234
$ws.Range("I6").Value = 0.04224437615225601
$ws.Range("J6").Value = 0.02856515650694651
$ws.Range("K6").Value = 1
$ws.Range("L6").Value = 0.5
$ws.Range("M6").Value = 0.358189
$ws.Range("N6").Value = 0.716378
$ws.Range("O6").Value = 0.052978466298774
$ws.Range("P6").Value = 0.03595390511601538
$ws.Range("Q6").Value = 0.122721640613
$ws.Range("R6").Value = 0.490886562452
$ws.Range("S6").Value = 0.002238042258295027
$ws.Range("T6").Value = 0.001027028926674884

# Row 7 (new)
$ws.Range("A7").Value = "MuSCs"
$ws.Range("B7").Value = "Cxcl13"
$ws.Range("C7").Value = "Cxcr3"
$ws.Range("D7").Value = "Inflammatory-Mac"
$ws.Range("E7").Value = 1
$ws.Range("F7").Value = 0.5
$ws.Range("G7").Value = 0.342617
$ws.Range("H7").Value = 0.685234
$ws.Range("I7").Value = 0.04224437615225601
$ws.Range("J7").Value = 0.02856515650694651
$ws.Range("K7").Value = 2
$ws.Range("L7").Value = 0.6666666666666666
$ws.Range("M7").Value = 0.574863
$ws.Range("N7").Value = 1.724589
$ws.Range("O7").Value = 0.08502595018806307
$ws.Range("P7").Value = 0.08655445766079338
$ws.Range("Q7").Value = 0.196957836471
$ws.Range("R7").Value = 1.181747018826
$ws.Range("S7").Value = 0.003591868222447519
$ws.Range("T7").Value = 0.002472441629454438

# Row 8 (new)
$ws.Range("A8").Value = "MuSCs"
$ws.Range("B8").Value = "Cxcl13"
$ws.Range("C8").Value = "Cxcr3"
$ws.Range("D8").Value = "Neutrophils"
$ws.Range("E8").Value = 1
$ws.Range("F8").Value = 0.5
$ws.Range("G8").Value = 0.342617
$ws.Range("H8").Value = 0.685234
$ws.Range("I8").Value = 0.04224437615225601
$ws.Range("J8").Value = 0.02856515650694651
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 5.590969
$ws.Range("N8").Value = 16.772907
$ws.Range("O8").Value = 0.8269404217996371
$ws.Range("P8").Value = 0.8418062905306279
$ws.Range("Q8").Value = 1.915561025873
$ws.Range("R8").Value = 11.493366155238
$ws.Range("S8").Value = 0.03493358223400912
$ws.Range("T8").Value = 0.02404632843753947

# Row 9 (new)
$ws.Range("A9").Value = "MuSCs"
$ws.Range("B9").Value = "Cxcl13"
$ws.Range("C9").Value = "Cxcr3"
$ws.Range("D9").Value = "Resolving-Mac"
$ws.Range("E9").Value = 1
$ws.Range("F9").Value = 0.5
$ws.Range("G9").Value = 0.342617
$ws.Range("H9").Value = 0.685234
$ws.Range("I9").Value = 0.04224437615225601
$ws.Range("J9").Value = 0.02856515650694651
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 0.237009
$ws.Range("N9").Value = 0.7110270000000001
$ws.Range("O9").Value = 0.0350551617135259
$ws.Range("P9").Value = 0.03568534669256324
$ws.Range("Q9").Value = 0.081203312553
$ws.Range("R9").Value = 0.487219875318
$ws.Range("S9").Value = 0.001480883437504352
$ws.Range("T9").Value = 0.001019357513277715

